# Update the "Förändrad" (Changed) date column (C) for rows 2 through 18
# from 45212 (2023-10-13) to 45221 (2023-10-22) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45212) {
        $cell.Value = 45221
    }
}
